# Applies the 2016-01 wordlist update:
#  - appends new vocabulary rows (B3:B10) to Arkusz1, backed by new shared strings
#  - gives the last two new rows (B9:B10) a distinct cell style (new font/xf)
#  - turns on iterative calculation (iterateDelta) for the workbook
#  - moves the active selection to B11 (just past the new data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    "Wir kommen gleich nach!",
    "verdorben",
    "ersticken",
    "Gefallen tun",
    "empören gegen",
    "hinrichten",
    "das Übel an der Wurzel packen",
    "Dolmetscher"
)

$startRow = 3
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newWords[$i]
}

# Give the last two rows (the idiom + the noun) their own cell style, same
# visual font but tracked as a distinct font/xf entry in styles.xml.
$styleRange = $ws.Range("B9:B10")
$styleRange.Font.Name = "Czcionka tekstu podstawowego"
$styleRange.Font.Family = 0

# Enable iterative calculation with a tighter convergence delta
# (Excel default maxChange is 0.001; the workbook now asks for 1E-4).
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# Move the selection past the newly-added data, mirroring the saved
# worksheet view (B11 is the next empty cell in column B).
$null = $ws.Range("B11").Select()
